$wb = $excel.ActiveWorkbook

# The old sheet carries baggage (merged cells, explicit column widths,
# leftover selection, etc.) that isn't exposed for in-place removal via
# the object model, so rebuild it from scratch: insert a fresh sheet,
# point it at the "Sheet1" name, and drop the original.
$oldSheet = $wb.Worksheets.Item("Sheet1")
$oldSheet.Name = "Sheet1_old"

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet1"

$wb.Worksheets.Item("Sheet1_old").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# New header row
$ws.Range("A1").Value = "Назва"
$ws.Range("B1").Value = "Значення"

# New data rows
$ws.Range("A2").Value = "А"
$ws.Range("B2").Value = 100

$ws.Range("A3").Value = "Б"
$ws.Range("B3").Value = 200

$ws.Range("A4").Value = "В"
$ws.Range("B4").Value = 300

# Header formatting: bold, bordered, centered horizontally, top-aligned vertically
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
